$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '36.767.54'
$ws.Range('E2').Value = '  +4.13%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.927.34'
$ws.Range('E3').Value = '  +2.55%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('B5').Value = 'BNB'
$ws.Range('C5').Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '250.38'
$ws.Range('E5').Value = '  +1.61%  '
$ws.Range('B6').Value = 'XRP'
$ws.Range('C6').Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.703'
$ws.Range('E6').Value = '  +2.86%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '44.42'
$ws.Range('E8').Value = '  +1.77%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '58.69'
$ws.Range('E9').Value = '  +9.19%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.370'
$ws.Range('E10').Value = '  +4.01%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0769'
$ws.Range('E11').Value = '  +4.10%  '
$ws.Range('E12').Value = '  +2.99%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '14.69'
$ws.Range('E13').Value = '  +8.56%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.820'
$ws.Range('E14').Value = '  +7.53%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.205.93'
$ws.Range('E15').Value = '  +2.49%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.15'
$ws.Range('E16').Value = '  +4.94%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.922.69'
$ws.Range('E17').Value = '  +2.12%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '36.770.25'
$ws.Range('E18').Value = '  +4.08%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '74.62'
$ws.Range('E19').Value = '  +2.57%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0868'
$ws.Range('E20').Value = '  +5.70%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '252.41'
$ws.Range('E21').Value = '  +3.47%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '13.46'
$ws.Range('E22').Value = '  +4.75%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.25'
$ws.Range('E23').Value = '  +5.82%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.68'
$ws.Range('E24').Value = '  +1.93%  '
$ws.Range('E25').Value = '  +0.01%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.22'
$ws.Range('E26').Value = '  +1.14%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '168.61'
$ws.Range('E27').Value = '  +1.69%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.90'
$ws.Range('E28').Value = '  +4.26%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '18.97'
$ws.Range('E29').Value = '  +3.83%  '
$ws.Range('E30').Value = '  +2.29%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.59'
$ws.Range('E31').Value = '  +6.78%  '
$ws.Range('B33').Value = 'InternetComputer(DFINITY)'
$ws.Range('C33').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.39'
$ws.Range('E33').Value = '  +5.87%  '
$ws.Range('B34').Value = 'WEMIXToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.96'
$ws.Range('E34').Value = '  -4.38%  '
$ws.Range('E35').Value = '  -0.09%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0876'
$ws.Range('E36').Value = '  +21.18%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.51'
$ws.Range('E37').Value = '  -11.70%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.900'
$ws.Range('E38').Value = '  +7.48%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '17.90'
$ws.Range('E39').Value = '  +50.05%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.06'
$ws.Range('E40').Value = '  +6.25%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '107.32'
$ws.Range('E41').Value = '  +11.73%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0229'
$ws.Range('E42').Value = '  +5.42%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '17.40'
$ws.Range('E43').Value = '  -1.83%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.12'
$ws.Range('E44').Value = '  +4.13%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.344.64'
$ws.Range('E45').Value = '  +3.22%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.58'
$ws.Range('E46').Value = '  +8.51%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.40'
$ws.Range('E47').Value = '  +1.46%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0817'
$ws.Range('E48').Value = '  +2.69%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.80'
$ws.Range('E49').Value = '  +2.55%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '6.45'
$ws.Range('E50').Value = '  +3.69%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '43.43'
$ws.Range('E51').Value = '  +3.27%  '
